$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record at the top of the series (row 14) and push the
# remaining historical rows down by one, which moves the former last row
# (38) into a brand new row 39.

$cols = @("D","J","K","L","M","O","P")
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot current values for rows 14..38 (the tracked columns) before
# overwriting anything.
$snapshot = @{}
for ($r = 14; $r -le 38; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Row 39 becomes a full copy of the old row 38 (every column).
foreach ($c in $allCols) {
    $ws.Range("$c" + "39").Value = $ws.Range("$c" + "38").Value()
}
$ws.Range("D39").NumberFormat = $ws.Range("D38").NumberFormat()

# Shift rows 15..38 down from the snapshot of rows 14..37 (process from the
# bottom up so we never read a value we've already overwritten).
for ($r = 38; $r -ge 15; $r--) {
    $src = $snapshot[$r - 1]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $src[$c]
    }
}

# Row 14 receives the brand-new weekly record.
$ws.Range("D14").Value = 44775
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("O14").Value = "Región de Coquimbo"
$ws.Range("P14").Value = 278
